$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$words = @(
    @("table", "Furniture"),
    @("flower", "Plant"),
    @("spider", "Animal"),
    @("kidney", "Organ"),
    @("moon", "Space"),
    @("fountain", "Water"),
    @("butterfly", "insect"),
    @("Umbrella", "Rain"),
    @("Necklace", "Jewelry"),
    @("Laptop", "Technology"),
    @("Guitar", "Music")
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 1
    $pair = $words[$i]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}
